# Prev_SO sheet: the treeview widget re-synced row 13 (cleared the stray
# blank placeholder row) and appended a new order, SO240307002, as row 15.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 was a leftover blank entry (empty inline-string cells in every
# column) - clear it out so the row carries no values.
for ($col = 1; $col -le 23; $col++) {
    $ws.Cells.Item(13, $col).Value = ""
}

# Row 15: new Sales Order entry appended by the treeview widget.
# Column A holds the date as literal text (not an Excel date serial), so
# force the cell to Text before writing it.
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "03/07/2024"
$ws.Range("B15").Value = "SO240307002"
$ws.Range("C15").Value = "NO"
$ws.Range("D15").Value = "ab"
$ws.Range("E15").Value = "(787)978-9777"
$ws.Range("F15").Value = "NO"
$ws.Range("G15").Value = "NO"
$ws.Range("H15").Value = "artist"
$ws.Range("I15").Value = "title"
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 87
$ws.Range("L15").Value = " "
$ws.Range("M15").Value = "AMA"
$ws.Range("N15").Value = "CD"
$ws.Range("O15").Value = "abake"
$ws.Range("P15").Value = "PICKUP"
